# Update cryptos list data (prices & volume change %) as per the
# "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that are plain text in the source
# workbook (e.g. "64.912.76" using '.' as a thousands separator, which is
# not a valid number). Assigning a numeric-looking string via .Value would
# otherwise get auto-coerced to a real number (and pick up a "Text" number
# format / style). Force General->Text->General so the stored cell stays a
# literal string with its original (default) style, matching the source.
function Set-TextValue {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

function Set-PctVol {
    param($row, $text)
    $ws.Range("E$row").Value = "  $text  "
}

# Row 2 - Bitcoin
Set-TextValue "D2" "64.802.62"
Set-PctVol 2 "-0.39%"

# Row 3 - Ethereum
Set-TextValue "D3" "3.513.50"
Set-PctVol 3 "-1.17%"

# Row 4 - TetherUSD
Set-PctVol 4 "+0.11%"

# Row 5 - BNB
Set-TextValue "D5" "595.89"
Set-PctVol 5 "-0.06%"

# Row 6 - Solana
Set-TextValue "D6" "133.15"
Set-PctVol 6 "-3.08%"

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.512.63"
Set-PctVol 7 "-1.11%"

# Row 8 - USDC
Set-PctVol 8 "+0.03%"

# Row 9 - XRP
Set-PctVol 9 "-0.63%"

# Row 10 - Dogecoin
Set-PctVol 10 "+0.22%"

# Row 11 - Toncoin
Set-PctVol 11 "+2.36%"

# Row 12 - Cardano
Set-TextValue "D12" "0.381"
Set-PctVol 12 "-1.32%"

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.115.14"
Set-PctVol 13 "-0.99%"

# Row 14 - Avalanche
Set-PctVol 14 "+0.60%"

# Row 15 - ShibaInu
Set-TextValue "D15" "0.0000180"
Set-PctVol 15 "-1.36%"

# Row 16 - TRON
Set-PctVol 16 "-0.01%"

# Row 17 - WrappedEther
Set-TextValue "D17" "3.515.21"
Set-PctVol 17 "-0.96%"

# Row 18 - WrappedBTC
Set-TextValue "D18" "64.868.63"
Set-PctVol 18 "+0.01%"

# Row 19 - Uniswap
Set-TextValue "D19" "10.04"
Set-PctVol 19 "-0.81%"

# Row 20 - Chainlink
Set-TextValue "D20" "14.28"
Set-PctVol 20 "-0.04%"

# Row 21 - Polkadot
Set-TextValue "D21" "5.66"
Set-PctVol 21 "-3.10%"

# Row 22 - BitcoinCash
Set-TextValue "D22" "390.31"
Set-PctVol 22 "-0.26%"

# Row 23 - Polygon
Set-PctVol 23 "-0.54%"

# Row 24 - WrappedeETH
Set-TextValue "D24" "3.657.10"
Set-PctVol 24 "-0.99%"

# Row 25 - Litecoin
Set-TextValue "D25" "73.88"
Set-PctVol 25 "-0.23%"

# Row 26 - Dai
Set-TextValue "D26" "0.999"
Set-PctVol 26 "-0.17%"

# Row 27 - PEPE
Set-TextValue "D27" "0.0000110"
Set-PctVol 27 "-3.98%"

# Row 28 - now RenderToken (was Fetch.AI)
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D28" "7.65"
Set-PctVol 28 "-1.45%"

# Row 29 - now Fetch.AI (was RenderToken)
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D29" "1.57"
Set-PctVol 29 "+11.34%"

# Row 30 - Binance-PegBSC-USD
Set-TextValue "D30" "0.998"
Set-PctVol 30 "-0.13%"

# Row 31 - PancakeSwap
Set-TextValue "D31" "2.27"
Set-PctVol 31 "-0.26%"

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "8.27"
Set-PctVol 32 "+0.07%"

# Row 33 - RenzoRestakedETH
Set-TextValue "D33" "3.518.88"
Set-PctVol 33 "-1.27%"

# Row 34 - EthereumClassic
Set-TextValue "D34" "24.11"
Set-PctVol 34 "+0.54%"

# Row 35 - USDe
Set-PctVol 35 "+0.01%"

# Row 36 - Kaspa
Set-TextValue "D36" "0.143"
Set-PctVol 36 "-0.45%"

# Row 37 - ImmutableX
Set-TextValue "D37" "1.56"
Set-PctVol 37 "+0.50%"

# Row 38 - now Monero (was NEARProtocol)
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D38" "168.57"
Set-PctVol 38 "-0.80%"

# Row 39 - now NEARProtocol (was Monero)
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "5.09"
Set-PctVol 39 "+1.69%"

# Row 40 - Aptos
Set-TextValue "D40" "6.81"
Set-PctVol 40 "-1.23%"

# Row 41 - Hedera
Set-TextValue "D41" "0.0814"
Set-PctVol 41 "+0.78%"

# Row 42 - Mantle
Set-TextValue "D42" "0.822"
Set-PctVol 42 "-0.28%"

# Row 43 - now OKB (was EnergySwap)
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "42.68"
Set-PctVol 43 "-0.72%"

# Row 44 - now EnergySwap (was ONDO)
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "25.60"
Set-PctVol 44 "-5.62%"

# Row 45 - now ONDO (was OKB)
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D45" "1.23"
Set-PctVol 45 "+0.50%"

# Row 46 - FirstDigitalUSD
Set-PctVol 46 "+0.14%"

# Row 47 - Filecoin
Set-TextValue "D47" "4.38"
Set-PctVol 47 "-1.48%"

# Row 48 - Stacks
Set-TextValue "D48" "1.64"
Set-PctVol 48 "-2.00%"

# Row 49 - Cosmos
Set-TextValue "D49" "6.87"
Set-PctVol 49 "-0.40%"

# Row 50 - Maker
Set-TextValue "D50" "2.374.36"
Set-PctVol 50 "-3.00%"

# Row 51 - VeChain
Set-TextValue "D51" "0.0267"
Set-PctVol 51 "+0.68%"
